$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '27.550.37'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '1.622.94'
$ws.Range('E3').Value = '  -1.44%  '
$ws.Range('D5').Value = '211.55'
$ws.Range('D6').Value = '0.527'
$ws.Range('E6').Value = '  -0.61%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('D8').Value = '23.20'
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').Value = '0.264'
$ws.Range('E9').Value = '  +2.15%  '
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('D11').Value = '0.0889'
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('D12').Value = '1.852.79'
$ws.Range('E12').Value = '  -1.47%  '
$ws.Range('D13').Value = '1.626.34'
$ws.Range('E13').Value = '  -1.24%  '
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('E15').Value = '  -1.82%  '
$ws.Range('D16').Value = '65.20'
$ws.Range('E16').Value = '  +0.78%  '
$ws.Range('D17').Value = '27.523.37'
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('D18').Value = '231.95'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('E20').Value = '  -0.70%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').Value = '10.30'
$ws.Range('E22').Value = '  +2.02%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '4.33'
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('E24').Value = '  +6.38%  '
$ws.Range('D25').Value = '150.13'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').Value = '6.87'
$ws.Range('E26').Value = '  -0.73%  '
$ws.Range('E27').Value = '  -0.44%  '
$ws.Range('D28').Value = '15.55'
$ws.Range('E28').Value = '  -0.60%  '
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').Value = '  -1.01%  '
$ws.Range('D31').Value = '0.0483'
$ws.Range('E31').Value = '  -0.77%  '
$ws.Range('E32').Value = '  -0.64%  '
$ws.Range('D33').Value = '1.474.09'
$ws.Range('E33').Value = '  +2.04%  '
$ws.Range('E34').Value = '  -2.12%  '
$ws.Range('E35').Value = '  -2.95%  '
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('E37').Value = '  +7.23%  '
$ws.Range('E38').Value = '  +0.40%  '
$ws.Range('D39').Value = '0.871'
$ws.Range('E39').Value = '  -0.96%  '
$ws.Range('D40').Value = '0.555'
$ws.Range('E40').Value = '  -2.21%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').Value = '1.02'
$ws.Range('E42').Value = '  -1.74%  '
$ws.Range('D43').Value = '67.69'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('B44').Value = 'mCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D44').Value = '2.46'
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').Value = '2.21'
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('D46').Value = '5.27'
$ws.Range('E46').Value = '  -5.52%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '1.763.88'
$ws.Range('E47').Value = '  -1.46%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '1.75'
$ws.Range('E48').Value = '  +0.72%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '87.29'
$ws.Range('E49').Value = '  +2.16%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0105'
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.101'
$ws.Range('E51').Value = '  +1.80%  '
